$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '29.572.75'
Set-TextValue 'E2' '  +4.32%  '

# Row 3
Set-TextValue 'D3' '1.602.55'
Set-TextValue 'E3' '  +3.45%  '

# Row 4
Set-TextValue 'E4' '  -0.46%  '

# Row 5
Set-TextValue 'D5' '213.66'
Set-TextValue 'E5' '  +1.85%  '

# Row 6
Set-TextValue 'D6' '0.515'
Set-TextValue 'E6' '  +7.04%  '

# Row 7
Set-TextValue 'E7' '  -0.43%  '

# Row 8
Set-TextValue 'D8' '26.80'
Set-TextValue 'E8' '  +12.21%  '

# Row 9
Set-TextValue 'E9' '  +3.53%  '

# Row 10
Set-TextValue 'E10' '  +2.95%  '

# Row 11
Set-TextValue 'E11' '  +2.87%  '

# Row 12
Set-TextValue 'D12' '1.831.89'
Set-TextValue 'E12' '  +3.51%  '

# Row 13
Set-TextValue 'D13' '1.605.56'
Set-TextValue 'E13' '  +4.01%  '

# Row 14
Set-TextValue 'D14' '29.587.59'
Set-TextValue 'E14' '  +4.48%  '

# Row 15
Set-TextValue 'E15' '  +3.80%  '

# Row 16
Set-TextValue 'E16' '  +3.59%  '

# Row 17
Set-TextValue 'D17' '63.55'
Set-TextValue 'E17' '  +4.41%  '

# Row 18
Set-TextValue 'D18' '243.71'
Set-TextValue 'E18' '  +6.96%  '

# Row 19
Set-TextValue 'D19' '7.60'
Set-TextValue 'E19' '  +3.45%  '

# Row 20
Set-TextValue 'D20' '0.0₃0696'
Set-TextValue 'E20' '  +2.95%  '

# Row 21
Set-TextValue 'E21' '  -0.50%  '

# Row 22
Set-TextValue 'E22' '  +3.84%  '

# Row 23
Set-TextValue 'E23' '  +3.96%  '

# Row 24
Set-TextValue 'E24' '  +4.35%  '

# Row 25
Set-TextValue 'D25' '155.33'
Set-TextValue 'E25' '  +2.54%  '

# Row 26
Set-TextValue 'D26' '15.36'
Set-TextValue 'E26' '  +4.34%  '

# Row 27
Set-TextValue 'E27' '  +5.75%  '

# Row 28
Set-TextValue 'D28' '6.39'
Set-TextValue 'E28' '  +2.57%  '

# Row 29
Set-TextValue 'D29' '0.997'
Set-TextValue 'E29' '  -0.32%  '

# Row 30
Set-TextValue 'D30' '0.0473'
Set-TextValue 'E30' '  +1.33%  '

# Row 31
Set-TextValue 'E31' '  +0.42%  '

# Row 32
Set-TextValue 'E32' '  +2.57%  '

# Row 33
Set-TextValue 'D33' '1.438.87'
Set-TextValue 'E33' '  +3.93%  '

# Row 34
Set-TextValue 'D34' '3.11'
Set-TextValue 'E34' '  +3.58%  '

# Row 35
Set-TextValue 'E35' '  -2.13%  '

# Row 36
Set-TextValue 'E36' '  +10.77%  '

# Row 37
Set-TextValue 'D37' '1.51'
Set-TextValue 'E37' '  +2.66%  '

# Row 38
Set-TextValue 'D38' '2.29'
Set-TextValue 'E38' '  -1.63%  '

# Row 39
Set-TextValue 'E39' '  +2.68%  '

# Row 40
Set-TextValue 'D40' '0.533'
Set-TextValue 'E40' '  +4.79%  '

# Row 41
Set-TextValue 'E41' '  +2.22%  '

# Row 42
Set-TextValue 'D42' '54.41'
Set-TextValue 'E42' '  +29.39%  '

# Row 43
Set-TextValue 'E43' '  +3.37%  '

# Row 44
Set-TextValue 'D44' '0.996'
Set-TextValue 'E44' '  -0.40%  '

# Row 45
Set-TextValue 'D45' '0.0467'
Set-TextValue 'E45' '  +3.06%  '

# Row 46
Set-TextValue 'D46' '65.81'
Set-TextValue 'E46' '  +6.31%  '

# Row 47
Set-TextValue 'E47' '  -0.16%  '

# Row 48
Set-TextValue 'D48' '1.742.55'
Set-TextValue 'E48' '  +3.66%  '

# Row 49
Set-TextValue 'D49' '86.40'
Set-TextValue 'E49' '  +0.79%  '

# Row 50
Set-TextValue 'D50' '0.838'
Set-TextValue 'E50' '  -3.07%  '

# Row 51
Set-TextValue 'E51' '  +0.90%  '
